$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (columns A-H)
$data = @(
    @("1328848", "https://aiesec.org/opportunity/global-talent/1328848", "Software Developer, machine vision (EU ONLY)", "Leuven, Belgium", "No", "5 applicants", "6 - 18 Months", "Heliovision"),
    @("1327778", "https://aiesec.org/opportunity/global-talent/1327778", "Digital Content & Stakeholder Engagement Intern", "Colombo, Sri Lanka", "No", "15 applicants", "6 - 18 Months", "Solutions Ground (Pvt) Ltd"),
    @("1327286", "https://aiesec.org/opportunity/global-talent/1327286", "Taxes & Internal Control", "Panamá, Provincia de Panamá, Panamá", "No", "58 applicants", "6 - 18 Months", "NESTLÉ"),
    @("1327281", "https://aiesec.org/opportunity/global-talent/1327281", "Purchasing Coordinator", "Bogotá, Colombia", "No", "37 applicants", "6 - 18 Months", "Microbiologia Y Genetica LTDA"),
    @("1327183", "https://aiesec.org/opportunity/global-talent/1327183", "ENGINEERING", "Çayırova, Kocaeli, Türkiye", "No", "52 applicants", "3 - 6 Months", "Eurotray Metal Elektrik San. Tic.: Ltd. Şti."),
    @("1325700", "https://aiesec.org/opportunity/global-talent/1325700", "Business Analyst and Executive Secretary", "Colombo, Sri Lanka", "No", "26 applicants", "3 - 6 Months", "Indian Kitchen PVT LTD"),
    @("1317170", "https://aiesec.org/opportunity/global-talent/1317170", "Guest Relations Officer", "Colombo, Sri Lanka", "No", "46 applicants", "3 - 6 Months", "Lanka Island Resorts Ltd"),
    @("1313206", "https://aiesec.org/opportunity/global-talent/1313206", "Digital Media Strategist", "Colombo, Sri Lanka", "No", "42 applicants", "9 - 12 Weeks", "Brand Corridor (Pvt) Ltd"),
    @("1280027", "https://aiesec.org/opportunity/global-talent/1280027", "Web Developer and Software in Drupart Company", "Gebze, Türkiye", "No", "680 applicants", "6 - 18 Months", "Drupart Dijital Çözümler ve Tic.Ltd.Şti.")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $val = $rowData[$c]
        if ($col -eq 1) {
            # Keep the opportunity ID column as text (not a number)
            $ws.Cells.Item($row, $col).Value = "'" + $val
        } else {
            $ws.Cells.Item($row, $col).Value = $val
        }
    }
}

# Remove rows 11 and 12 (previously had 12 rows, now only 10 remain)
$ws.Rows.Item(11).Delete() | Out-Null
$ws.Rows.Item(11).Delete() | Out-Null

# Update column widths
# Note: the headless ColumnWidth setter adds a constant +0.8333333 (5/6)
# pixel-padding offset relative to the stored OOXML "width" attribute, so we
# subtract it here to land exactly on the target stored widths (50/38/17/16/47).
$widthOffset = 0.8333333
$ws.Columns.Item(3).ColumnWidth = 50 - $widthOffset
$ws.Columns.Item(4).ColumnWidth = 38 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 17 - $widthOffset
$ws.Columns.Item(7).ColumnWidth = 16 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 47 - $widthOffset
